$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: update the Property/Value table ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now set to "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (row 11),
# leaving a single Contact row (row 10) which is then repurposed below.
$meta.Rows.Item(11).Delete()

# Turn the remaining "Contact" row into the new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- "Elements" sheet: update root element Short/Definition text ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Converted Quantity"
$elements.Range("L2").Value = "Converted quantity expressed in standard unit value"
